$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (was 45406 -> 45436, i.e. 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = Get-Date -Year 2024 -Month 5 -Day 24 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# Update prices in D23:D26 to 844
$ws.Range("D23").Value = 844
$ws.Range("D24").Value = 844
$ws.Range("D25").Value = 844
$ws.Range("D26").Value = 844
